$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (A8:C8) previously held the "z16" record, which the column guesser
# choked on (NPE). Replace it with the last record (previously row 18,
# "z26"/"z36"/"z46") and drop the now-duplicated trailing row.
$ws.Range("A8:C8").Value2 = $ws.Range("A18:C18").Value2

$ws.Rows("18:18").Delete()

$ws.Range("A8:XFD8").Select()
